$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header cell A1 from "ID" to "Id"
$ws.Range("A1").Value = "Id"

# Move the active selection to P22 (mirrors the final selection in the saved file)
$ws.Range("P22").Select()
